$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" header columns (A1:J1) to "_FV2410" and the "_new"
# header columns (L1:U1) to "_FV2504". Column K1 ("diff") is unchanged.
$oldHeaders = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
$newHeaders = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# Turn the header/data range into an Excel Table ("Table1").
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U64"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split at row 1, frozen).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
